$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date values for rows 2-5 from 2023-09-16 (45185) to 2023-10-05 (45204),
# preserving existing cell formatting. Use the raw Excel serial number so no time component is added.
$ws.Range("C2").Value = 45204
$ws.Range("C3").Value = 45204
$ws.Range("C4").Value = 45204
$ws.Range("C5").Value = 45204
